# Applies scheduled-runner price/profit refresh values to the Zeromus_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 300.2857
$ws.Range("I31").Value = 300.2857
$ws.Range("K31").Value = 900.8571000000001
$ws.Range("M31").Value = -670.8571000000001

# Row 132
$ws.Range("H132").Value = 4479.2544
$ws.Range("I132").Value = 4768.975
$ws.Range("J132").Value = 3706.6667
$ws.Range("K132").Value = 14306.925
$ws.Range("L132").Value = 11120.0001
$ws.Range("M132").Value = -11776.925
$ws.Range("N132").Value = -16180.0001

# Row 138
$ws.Range("H138").Value = 3039.2783
$ws.Range("I138").Value = 1873.4849
$ws.Range("J138").Value = 3640.3906
$ws.Range("K138").Value = 5620.4547
$ws.Range("L138").Value = 10921.1718
$ws.Range("M138").Value = -480.4547000000002
$ws.Range("N138").Value = -21201.1718

$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 7684.3335
$ws.Range("J37").Value = 17969
$ws.Range("L37").Value = 17969
$ws.Range("N37").Value = -18515

# Row 61
$ws.Range("H61").Value = 1808.0344
$ws.Range("I61").Value = 1476.65
$ws.Range("J61").Value = 2544.4443
$ws.Range("K61").Value = 1476.65
$ws.Range("L61").Value = 2544.4443
$ws.Range("M61").Value = -1264.65
$ws.Range("N61").Value = -2968.4443

# Row 136
$ws.Range("H136").Value = 1808.0344
$ws.Range("I136").Value = 1476.65
$ws.Range("J136").Value = 2544.4443
$ws.Range("K136").Value = 4429.950000000001
$ws.Range("L136").Value = 7633.3329
$ws.Range("M136").Value = -1879.950000000001
$ws.Range("N136").Value = -12733.3329

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 239048.22
$ws.Range("I134").Value = 271008.78
$ws.Range("J134").Value = 2540
$ws.Range("K134").Value = 813026.3400000001
$ws.Range("L134").Value = 7620
$ws.Range("M134").Value = -810491.3400000001
$ws.Range("N134").Value = -12690

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 1130
$ws.Range("I107").Value = 1056.9231
$ws.Range("J107").Value = 1320
$ws.Range("K107").Value = 1056.9231
$ws.Range("L107").Value = 1320
$ws.Range("M107").Value = 863.0769
$ws.Range("N107").Value = -5160

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 82518.125
$ws.Range("I2").Value = 198012.8
$ws.Range("K2").Value = 1188076.8
$ws.Range("M2").Value = -1187963.8

# Row 5
$ws.Range("H5").Value = 898.51514
$ws.Range("J5").Value = 1388.7693
$ws.Range("L5").Value = 4166.3079
$ws.Range("N5").Value = -4390.3079

# Row 22
$ws.Range("H22").Value = 4166.6665
$ws.Range("I22").Value = 5500.5
$ws.Range("J22").Value = 3899.9
$ws.Range("K22").Value = 16501.5
$ws.Range("L22").Value = 11699.7
$ws.Range("M22").Value = -16332.5
$ws.Range("N22").Value = -12037.7

# Row 27
$ws.Range("H27").Value = 4166.6665
$ws.Range("I27").Value = 5500.5
$ws.Range("J27").Value = 3899.9
$ws.Range("K27").Value = 16501.5
$ws.Range("L27").Value = 11699.7
$ws.Range("M27").Value = -16399.5
$ws.Range("N27").Value = -11903.7

# Row 33
$ws.Range("H33").Value = 6881.7334
$ws.Range("I33").Value = 204.33333
$ws.Range("J33").Value = 11333.333
$ws.Range("K33").Value = 1225.99998
$ws.Range("L33").Value = 67999.99800000001
$ws.Range("M33").Value = -942.9999800000001
$ws.Range("N33").Value = -68565.99800000001

# Row 44
$ws.Range("H44").Value = 5785.533
$ws.Range("I44").Value = 1494.75
$ws.Range("J44").Value = 7345.8184
$ws.Range("K44").Value = 4484.25
$ws.Range("L44").Value = 22037.4552
$ws.Range("M44").Value = -4086.25
$ws.Range("N44").Value = -22833.4552

# Row 58
$ws.Range("H58").Value = 3741.6667
$ws.Range("I58").Value = 1850
$ws.Range("J58").Value = 3978.125
$ws.Range("K58").Value = 5550
$ws.Range("L58").Value = 11934.375
$ws.Range("M58").Value = -5422
$ws.Range("N58").Value = -12190.375

# Row 64
$ws.Range("H64").Value = 5108.8184
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 5256.857
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 15770.571
$ws.Range("M64").Value = -5730
$ws.Range("N64").Value = -16310.571

# Row 67
$ws.Range("H67").Value = 5108.8184
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 5256.857
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 15770.571
$ws.Range("M67").Value = -5064
$ws.Range("N67").Value = -17642.571

# Row 68
$ws.Range("H68").Value = 1560.3
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 1700.375
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 5101.125
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -6723.125

# Row 70
$ws.Range("H70").Value = 5179
$ws.Range("I70").Value = 1861.2
$ws.Range("J70").Value = 6052.1055
$ws.Range("K70").Value = 5583.6
$ws.Range("L70").Value = 18156.3165
$ws.Range("M70").Value = -5268.6
$ws.Range("N70").Value = -18786.3165

# Row 71
$ws.Range("H71").Value = 1560.3
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 1700.375
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 15303.375
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -23415.375

# Row 73
$ws.Range("H73").Value = 5179
$ws.Range("I73").Value = 1861.2
$ws.Range("J73").Value = 6052.1055
$ws.Range("K73").Value = 5583.6
$ws.Range("L73").Value = 18156.3165
$ws.Range("M73").Value = -4491.6
$ws.Range("N73").Value = -20340.3165

# Row 76
$ws.Range("H76").Value = 5597.9443
$ws.Range("I76").Value = 2093.8333
$ws.Range("J76").Value = 7350
$ws.Range("K76").Value = 6281.499899999999
$ws.Range("L76").Value = 22050
$ws.Range("M76").Value = -5898.499899999999
$ws.Range("N76").Value = -22816

# Row 79
$ws.Range("H79").Value = 5597.9443
$ws.Range("I79").Value = 2093.8333
$ws.Range("J79").Value = 7350
$ws.Range("K79").Value = 6281.499899999999
$ws.Range("L79").Value = 22050
$ws.Range("M79").Value = -4955.499899999999
$ws.Range("N79").Value = -24702

# Row 92
$ws.Range("H92").Value = 554.2
$ws.Range("I92").Value = 491.7143
$ws.Range("K92").Value = 1475.1429
$ws.Range("M92").Value = -227.1428999999998

# Row 97
$ws.Range("H97").Value = 1250.75
$ws.Range("I97").Value = 1001
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 3003
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -2507
$ws.Range("N97").Value = -6992

# Row 103
$ws.Range("H103").Value = 1366.4667
$ws.Range("I103").Value = 302.25
$ws.Range("K103").Value = 906.75
$ws.Range("M103").Value = -27.75

# Row 106
$ws.Range("H106").Value = 5802.364
$ws.Range("J106").Value = 5802.364
$ws.Range("L106").Value = 17407.092
$ws.Range("N106").Value = -19299.092

# Row 112
$ws.Range("H112").Value = 5481.4688
$ws.Range("I112").Value = 4509
$ws.Range("J112").Value = 5582.069
$ws.Range("K112").Value = 13527
$ws.Range("L112").Value = 16746.207
$ws.Range("M112").Value = -12419
$ws.Range("N112").Value = -18962.207

# Row 116
$ws.Range("H116").Value = 4796.3335
$ws.Range("I116").Value = 2682.3333
$ws.Range("J116").Value = 6910.3335
$ws.Range("K116").Value = 8046.999899999999
$ws.Range("L116").Value = 20731.0005
$ws.Range("M116").Value = -4604.999899999999
$ws.Range("N116").Value = -27615.0005

# Row 121
$ws.Range("H121").Value = 1370.3077
$ws.Range("I121").Value = 392.23077
$ws.Range("J121").Value = 1859.3462
$ws.Range("K121").Value = 1176.69231
$ws.Range("L121").Value = 5578.0386
$ws.Range("M121").Value = 133.3076900000001
$ws.Range("N121").Value = -8198.0386

# Row 132
$ws.Range("H132").Value = 1227.381
$ws.Range("I132").Value = 637.5
$ws.Range("J132").Value = 1763.6364
$ws.Range("K132").Value = 5737.5
$ws.Range("L132").Value = 15872.7276
$ws.Range("M132").Value = -3207.5
$ws.Range("N132").Value = -20932.7276

# Row 135
$ws.Range("H135").Value = 898.51514
$ws.Range("J135").Value = 1388.7693
$ws.Range("L135").Value = 12498.9237
$ws.Range("N135").Value = -17568.9237

# Row 140
$ws.Range("H140").Value = 1493.475
$ws.Range("I140").Value = 1268.1111
$ws.Range("K140").Value = 3804.3333
$ws.Range("M140").Value = 1375.6667

$ws = $wb.Worksheets.Item("LTW")
# Row 32
$ws.Range("H32").Value = 8507.5
$ws.Range("J32").Value = 8507.5
$ws.Range("L32").Value = 8507.5
$ws.Range("N32").Value = -9141.5

# Row 61
$ws.Range("H61").Value = 2142.8948
$ws.Range("I61").Value = 1967.5
$ws.Range("J61").Value = 2443.5715
$ws.Range("K61").Value = 1967.5
$ws.Range("L61").Value = 2443.5715
$ws.Range("M61").Value = -1765.5
$ws.Range("N61").Value = -2847.5715

# Row 113
$ws.Range("H113").Value = 2142.8948
$ws.Range("I113").Value = 1967.5
$ws.Range("J113").Value = 2443.5715
$ws.Range("K113").Value = 1967.5
$ws.Range("L113").Value = 2443.5715
$ws.Range("M113").Value = 202.5
$ws.Range("N113").Value = -6783.5715

# Row 132
$ws.Range("H132").Value = 13897371
$ws.Range("I132").Value = 19241300
$ws.Range("J132").Value = 3152
$ws.Range("K132").Value = 57723900
$ws.Range("L132").Value = 9456
$ws.Range("M132").Value = -57721370
$ws.Range("N132").Value = -14516

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
